$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.384.29"
$ws.Range("E2").Value = "'  +0.00%  "
$ws.Range("D3").Value = "'1.571.02"
$ws.Range("E3").Value = "'  -0.11%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E5").Value = "'  +0.11%  "
$ws.Range("D6").Value = "'291.81"
$ws.Range("E6").Value = "'  +0.58%  "
$ws.Range("D7").Value = "'0.3765"
$ws.Range("E7").Value = "'  +2.29%  "
$ws.Range("D8").Value = "'49.71"
$ws.Range("E8").Value = "'  +0.63%  "
$ws.Range("E9").Value = "'  +0.58%  "
$ws.Range("D10").Value = "'0.07618"
$ws.Range("E10").Value = "'  -0.10%  "
$ws.Range("D11").Value = "'1.145"
$ws.Range("E11").Value = "'  -1.84%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "'  +0.13%  "
$ws.Range("D13").Value = "'21.12"
$ws.Range("E13").Value = "'  -1.30%  "
$ws.Range("D14").Value = "'6.005"
$ws.Range("E14").Value = "'  -0.96%  "
$ws.Range("D15").Value = "'6.968"
$ws.Range("E15").Value = "'  +0.70%  "
$ws.Range("D16").Value = "'1.568.22"
$ws.Range("E16").Value = "'  +0.08%  "
$ws.Range("E17").Value = "'  -0.19%  "
$ws.Range("E18").Value = "'  +0.56%  "
$ws.Range("D19").Value = "'0.06746"
$ws.Range("E19").Value = "'  -0.07%  "
$ws.Range("E20").Value = "'  +0.10%  "
$ws.Range("D21").Value = "'16.71"
$ws.Range("E21").Value = "'  +1.07%  "
$ws.Range("D22").Value = "'6.193"
$ws.Range("E22").Value = "'  -0.61%  "
$ws.Range("E23").Value = "'  -0.28%  "
$ws.Range("D24").Value = "'22.389.69"
$ws.Range("E24").Value = "'  -0.03%  "
$ws.Range("D25").Value = "'2.387"
$ws.Range("E25").Value = "'  +0.36%  "
$ws.Range("D26").Value = "'2.690"
$ws.Range("E26").Value = "'  -7.89%  "
$ws.Range("D27").Value = "'20.13"
$ws.Range("E27").Value = "'  +0.52%  "
$ws.Range("D28").Value = "'147.57"
$ws.Range("E28").Value = "'  +0.91%  "
$ws.Range("D29").Value = "'5.039"
$ws.Range("E29").Value = "'  +1.32%  "
$ws.Range("D30").Value = "'126.52"
$ws.Range("E30").Value = "'  +0.58%  "
$ws.Range("D31").Value = "'1.746.87"
$ws.Range("E31").Value = "'  +0.28%  "
$ws.Range("E32").Value = "'  +0.70%  "
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.101"
$ws.Range("E33").Value = "'  -2.77%  "
$ws.Range("B34").Value = "'ImmutableX"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.9977"
$ws.Range("E34").Value = "'  -4.73%  "
$ws.Range("D35").Value = "'10.10"
$ws.Range("E35").Value = "'  -1.30%  "
$ws.Range("D36").Value = "'0.08498"
$ws.Range("E36").Value = "'  -0.33%  "
$ws.Range("D37").Value = "'0.02535"
$ws.Range("E37").Value = "'  -0.06%  "
$ws.Range("D38").Value = "'1.390"
$ws.Range("E38").Value = "'  +10.73%  "
$ws.Range("E39").Value = "'  -1.00%  "
$ws.Range("E40").Value = "'  -0.56%  "
$ws.Range("D41").Value = "'5.402"
$ws.Range("E41").Value = "'  -2.62%  "
$ws.Range("D42").Value = "'11.40"
$ws.Range("E42").Value = "'  -2.81%  "
$ws.Range("D43").Value = "'0.6339"
$ws.Range("E43").Value = "'  -0.38%  "
$ws.Range("E44").Value = "'  +0.18%  "
$ws.Range("D45").Value = "'14.01"
$ws.Range("E45").Value = "'  -2.71%  "
$ws.Range("D46").Value = "'3.807"
$ws.Range("E46").Value = "'  +1.61%  "
$ws.Range("D47").Value = "'0.5947"
$ws.Range("E47").Value = "'  -0.72%  "
$ws.Range("B48").Value = "'EOS"
$ws.Range("C48").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "'1.284"
$ws.Range("E48").Value = "'  +1.50%  "
$ws.Range("B49").Value = "'NEARProtocol"
$ws.Range("C49").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.084"
$ws.Range("E49").Value = "'  -1.53%  "
$ws.Range("D50").Value = "'124.44"
$ws.Range("E50").Value = "'  +0.40%  "
$ws.Range("D51").Value = "'0.07318"
